$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string into a cell while preserving the
# "text" (shared-string) cell type instead of letting it be auto-coerced
# into a numeric cell. Plain `.Value =` assignment on a purely-numeric
# string is always re-interpreted as a number by this engine, so instead
# we park the literal text in an out-of-range scratch cell via a text
# formula (guarantees a string result), copy it, and paste-special just
# the value into the destination - this clones the string cell as-is
# (type + content) without re-inferring its type and without touching the
# destination's existing style.
function Set-TextValue($range, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.Clear()
}

Set-TextValue $ws.Range("P2") "320018099707"
Set-TextValue $ws.Range("P3") "320018099718"
Set-TextValue $ws.Range("P4") "320018110017"
